$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Block 1: columns B-F, rows 2-25
$data1 = New-Object 'object[,]' 24,5
$data1[0,0] = 1.02
$data1[0,1] = 1.026490916865734
$data1[0,2] = 1.030263256528299
$data1[0,3] = 1.035281354641555
$data1[0,4] = 1.043026249814855
$data1[1,0] = 1.02
$data1[1,1] = 1.027842156397847
$data1[1,2] = 1.030606612076188
$data1[1,3] = 1.03652365096261
$data1[1,4] = 1.044448464980237
$data1[2,0] = 1.02
$data1[2,1] = 1.028716173043294
$data1[2,2] = 1.030828240658591
$data1[2,3] = 1.037327489195902
$data1[2,4] = 1.045369012804245
$data1[3,0] = 1.02
$data1[3,1] = 1.02908353615165
$data1[3,2] = 1.030921280219724
$data1[3,3] = 1.037665423616462
$data1[3,4] = 1.045756082122717
$data1[4,0] = 1.02
$data1[4,1] = 1.029145213854658
$data1[4,2] = 1.030936894093018
$data1[4,3] = 1.037722164438539
$data1[4,4] = 1.045821077047701
$data1[5,0] = 1.02
$data1[5,1] = 1.028721082050322
$data1[5,2] = 1.030829484383945
$data1[5,3] = 1.037332004684755
$data1[5,4] = 1.045374184559297
$data1[6,0] = 1.02
$data1[6,1] = 1.026947644726116
$data1[6,2] = 1.030379406564799
$data1[6,3] = 1.03570119870601
$data1[6,4] = 1.043506838681758
$data1[7,0] = 1.02
$data1[7,1] = 1.023819926531252
$data1[7,2] = 1.029582245551909
$data1[7,3] = 1.032827278913259
$data1[7,4] = 1.04021829376622
$data1[8,0] = 1.02
$data1[8,1] = 1.021732687406677
$data1[8,2] = 1.029048220613477
$data1[8,3] = 1.030910964937144
$data1[8,4] = 1.038026988218159
$data1[9,0] = 1.02
$data1[9,1] = 1.020828330211808
$data1[9,2] = 1.028816400681362
$data1[9,3] = 1.030081045017767
$data1[9,4] = 1.037078321730686
$data1[10,0] = 1.02
$data1[10,1] = 1.020492320724775
$data1[10,2] = 1.028730207336787
$data1[10,3] = 1.029772750100411
$data1[10,4] = 1.036725967244162
$data1[11,0] = 1.02
$data1[11,1] = 1.020564400104801
$data1[11,2] = 1.028748699911505
$data1[11,3] = 1.029838881637243
$data1[11,4] = 1.036801547541351
$data1[12,0] = 1.02
$data1[12,1] = 1.020800557450148
$data1[12,2] = 1.028809277633772
$data1[12,3] = 1.030055561812355
$data1[12,4] = 1.037049195569475
$data1[13,0] = 1.02
$data1[13,1] = 1.020946049567315
$data1[13,2] = 1.028846590368057
$data1[13,3] = 1.030189062087951
$data1[13,4] = 1.037201782521287
$data1[14,0] = 1.02
$data1[14,1] = 1.021792694161462
$data1[14,2] = 1.029063593663635
$data1[14,3] = 1.030966040613384
$data1[14,4] = 1.038089951519511
$data1[15,0] = 1.02
$data1[15,1] = 1.022323615771931
$data1[15,2] = 1.029199559676678
$data1[15,3] = 1.031453377401597
$data1[15,4] = 1.038647121976612
$data1[16,0] = 1.02
$data1[16,1] = 1.022633238848842
$data1[16,2] = 1.029278809869972
$data1[16,3] = 1.03173761916036
$data1[16,4] = 1.038972128104004
$data1[17,0] = 1.02
$data1[17,1] = 1.02273880323527
$data1[17,2] = 1.029305822478735
$data1[17,3] = 1.031834536082859
$data1[17,4] = 1.039082950004158
$data1[18,0] = 1.02
$data1[18,1] = 1.022266658587825
$data1[18,2] = 1.029184977633886
$data1[18,3] = 1.031401092208852
$data1[18,4] = 1.038587341033059
$data1[19,0] = 1.02
$data1[19,1] = 1.020731017553108
$data1[19,2] = 1.02879144133297
$data1[19,3] = 1.029991755644657
$data1[19,4] = 1.036976268817576
$data1[20,0] = 1.02
$data1[20,1] = 1.019764970093797
$data1[20,2] = 1.0285435183423
$data1[20,3] = 1.029105498996336
$data1[20,4] = 1.035963451356113
$data1[21,0] = 1.02
$data1[21,1] = 1.020277142110416
$data1[21,2] = 1.028674992704237
$data1[21,3] = 1.029575336107435
$data1[21,4] = 1.036500354652359
$data1[22,0] = 1.02
$data1[22,1] = 1.022292395255737
$data1[22,2] = 1.029191566806252
$data1[22,3] = 1.031424717678507
$data1[22,4] = 1.038614353412237
$data1[23,0] = 1.02
$data1[23,1] = 1.024628864977666
$data1[23,2] = 1.029788797171083
$data1[23,3] = 1.033570305248743
$data1[23,4] = 1.041068256048337

$ws.Range("B2:F25").Value = $data1

# Block 2: columns I-N, rows 2-25
$data2 = New-Object 'object[,]' 24,6
$data2[0,0] = 1.02955528997503
$data2[0,1] = 1.031654064059758
$data2[0,2] = 1.033074774661431
$data2[0,3] = 1.03807840888099
$data2[0,4] = 1.045801277101354
$data2[0,5] = 1.033119131163119
$data2[1,0] = 1.029546965970127
$data2[1,1] = 1.032643018763664
$data2[1,2] = 1.033227698754684
$data2[1,3] = 1.039128916684577
$data2[1,4] = 1.047032844707646
$data2[1,5] = 1.034109490296118
$data2[2,0] = 1.029539753632458
$data2[2,1] = 1.033282201577143
$data2[2,2] = 1.0333253549493
$data2[2,3] = 1.039808101438636
$data2[2,4] = 1.047829512859632
$data2[2,5] = 1.034749580822508
$data2[3,0] = 1.029536283262584
$data2[3,1] = 1.033550740576031
$data2[3,2] = 1.033366097555123
$data2[3,3] = 1.040093498731202
$data2[3,4] = 1.048164378682872
$data2[3,5] = 1.035018501177576
$data2[4,0] = 1.029535674837201
$data2[4,1] = 1.033595819370691
$data2[4,2] = 1.033372920065691
$data2[4,3] = 1.040141410541831
$data2[4,4] = 1.048220601097659
$data2[4,5] = 1.035063643989295
$data2[5,0] = 1.029539708984861
$data2[5,1] = 1.033285790486255
$data2[5,2] = 1.033325900582478
$data2[5,3] = 1.039811915446027
$data2[5,4] = 1.047833987558758
$data2[5,5] = 1.034753174828282
$data2[6,0] = 1.02955285428455
$data2[6,1] = 1.031988440333347
$data2[6,2] = 1.033126723317437
$data2[6,3] = 1.038433551539637
$data2[6,4] = 1.046217542688354
$data2[6,5] = 1.033453982289349
$data2[7,0] = 1.029562081003761
$data2[7,1] = 1.029696565096052
$data2[7,2] = 1.032765901988229
$data2[7,3] = 1.036000241370182
$data2[7,4] = 1.043367172052731
$data2[7,5] = 1.031158852326335
$data2[8,0] = 1.029558927936093
$data2[8,1] = 1.028164574574658
$data2[8,2] = 1.032518842347417
$data2[8,3] = 1.034374845117876
$data2[8,4] = 1.041465360200861
$data2[8,5] = 1.029624686202679
$data2[9,0] = 1.029555369739793
$data2[9,1] = 1.027500197155752
$data2[9,2] = 1.032410340936442
$data2[9,3] = 1.033670230717229
$data2[9,4] = 1.040641432365732
$data2[9,5] = 1.028959365291633
$data2[10,0] = 1.02955371971951
$data2[10,1] = 1.027253261528089
$data2[10,2] = 1.032369811708957
$data2[10,3] = 1.033408380493541
$data2[10,4] = 1.040335319692193
$data2[10,5] = 1.028712078987055
$data2[11,0] = 1.029554088495152
$data2[11,1] = 1.027306237190017
$data2[11,2] = 1.032378515607504
$data2[11,3] = 1.033464553975997
$data2[11,4] = 1.040400985102832
$data2[11,5] = 1.028765129880499
$data2[12,0] = 1.029555240040407
$data2[12,1] = 1.027479788589222
$data2[12,2] = 1.032406995397466
$data2[12,3] = 1.033648588669355
$data2[12,4] = 1.04061613043405
$data2[12,5] = 1.028938927742598
$data2[13,0] = 1.029555906069241
$data2[13,1] = 1.027586698515009
$data2[13,2] = 1.032424512714591
$data2[13,3] = 1.033761961831765
$data2[13,4] = 1.040748679270769
$data2[13,5] = 1.029045989492742
$data2[14,0] = 1.029559117975007
$data2[14,1] = 1.028208645405253
$data2[14,2] = 1.032526011298701
$data2[14,3] = 1.034421590726031
$data2[14,4] = 1.041520032120772
$data2[14,5] = 1.029668819618909
$data2[15,0] = 1.029560546404202
$data2[15,1] = 1.028598501782261
$data2[15,2] = 1.032589272239751
$data2[15,3] = 1.034835139739719
$data2[15,4] = 1.042003762906309
$data2[15,5] = 1.030059229636689
$data2[16,0] = 1.029561167926103
$data2[16,1] = 1.02882580069113
$data2[16,2] = 1.032626024211051
$data2[16,3] = 1.035076278268112
$data2[16,4] = 1.042285873631966
$data2[16,5] = 1.03028685133608
$data2[17,0] = 1.029561343914426
$data2[17,1] = 1.02890328728206
$data2[17,2] = 1.032638530682611
$data2[17,3] = 1.035158487224475
$data2[17,4] = 1.042382059235768
$data2[17,5] = 1.030364447966857
$data2[18,0] = 1.029560415031944
$data2[18,1] = 1.02855668401827
$data2[18,2] = 1.032582500137955
$data2[18,3] = 1.034790777901108
$data2[18,4] = 1.041951867472192
$data2[18,5] = 1.030017352486677
$data2[19,0] = 1.029554909994162
$data2[19,1] = 1.02742868637165
$data2[19,2] = 1.032398615060048
$data2[19,3] = 1.033594398522078
$data2[19,4] = 1.040552777442377
$data2[19,5] = 1.028887752954016
$data2[20,0] = 1.029549549582982
$data2[20,1] = 1.026718563481539
$data2[20,2] = 1.032281686918611
$data2[20,3] = 1.032841461854167
$data2[20,4] = 1.039672711468772
$data2[20,5] = 1.028176621607948
$data2[21,0] = 1.029552570889732
$data2[21,1] = 1.027095100184892
$data2[21,2] = 1.032343796520178
$data2[21,3] = 1.033240677867888
$data2[21,4] = 1.040139290747095
$data2[21,5] = 1.028553693036609
$data2[22,0] = 1.029560475047559
$data2[22,1] = 1.028575579967526
$data2[22,2] = 1.032585560613934
$data2[22,3] = 1.034810823346057
$data2[22,4] = 1.04197531691229
$data2[22,5] = 1.030036275270348
$data2[23,0] = 1.029561340424493
$data2[23,1] = 1.030289773934296
$data2[23,2] = 1.032860338570589
$data2[23,3] = 1.036629858356135
$data2[23,4] = 1.044104320226557
$data2[23,5] = 1.031752903589171

$ws.Range("I2:N25").Value = $data2

Write-Host "case with 380 kV done"
